$wb = $excel.ActiveWorkbook
$tracksSheet = $wb.Worksheets.Item("Tracks")
$newSheet = $wb.Worksheets.Add($null, $tracksSheet)
$newSheet.Name = "CoursesTracks"

$newSheet.Cells.Item(1,1).Value = "track_id"
$newSheet.Cells.Item(1,2).Value = "course_number"

$trackIds = @(1,1,1,1,1,1,2,2,2,2,2,2,2,2,2,2,2,3,3,3,3,3,3,3,3,3,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,5,5,5,5,5,5,5,5,5,5,5,5,5,6,6,6,6,6,6,6,6,6,6,6,6,7,7,7,7,7,7,7,7,7,7,7,7)
$courseNums = @(302,303,307,340,442,441,302,307,312,313,340,341,342,344,441,442,443,302,303,318,344,361,345,400,405,462,300,320,321,301,316,317,400,401,405,402,403,420,421,440,302,303,318,428,450,300,301,302,303,307,316,317,345,360,400,401,405,341,307,308,312,340,322,302,414,323,441,344,442,316,307,340,301,317,344,302,303,320,321,441,402,403)

for ($i = 0; $i -lt $trackIds.Length; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $trackIds[$i]
    $newSheet.Cells.Item($row, 2).Value = $courseNums[$i]
}

$tracksSheet.Cells.Item(2,2).Value = "Mathematics of Information"
$tracksSheet.Cells.Item(3,1).Value = 2
$tracksSheet.Cells.Item(3,2).Value = "Mathematics of Discrete Algorithms"
$tracksSheet.Cells.Item(4,1).Value = 3
$tracksSheet.Cells.Item(4,2).Value = "Mathematical Biology"
$tracksSheet.Cells.Item(5,1).Value = 4
$tracksSheet.Cells.Item(5,2).Value = "Applied Track for Graduate School"
$tracksSheet.Cells.Item(6,1).Value = 5
$tracksSheet.Cells.Item(6,2).Value = "Mathematical Modelling and Copmutation"
$tracksSheet.Cells.Item(7,1).Value = 6
$tracksSheet.Cells.Item(7,2).Value = "Education"
$tracksSheet.Cells.Item(8,1).Value = 7
$tracksSheet.Cells.Item(8,2).Value = "Mathematical Optimization"
